$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    0.99999999005975926,
    0.9952069477150558,
    0.97889967873030748,
    0.97271255584430361,
    0.96705421883929332,
    0.95513517803678538,
    0.95371257868841397,
    0.95272132128839671,
    0.95483878971552172,
    0.95794113594119357,
    0.95856702585381603,
    0.96001606027942277,
    0.96884229233182306,
    0.97360902998503884,
    0.97101826090086796,
    0.968512258872122,
    0.96480484336172723,
    0.96369597431367793,
    0.99241918101576765,
    0.96950304950522836,
    0.96663326554822493,
    0.96292016323852936,
    0.97186456699314083,
    0.95884335699345624,
    0.95238631847498922,
    0.95343032062075184,
    0.95130390812995658,
    0.94436403292867621,
    0.93989797923571083,
    0.93802360984412669,
    0.94464376563813746,
    0.9426166117015764,
    0.94195139284440588
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $values[$i]
}
